$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 ("bl_high" / "NUMERIC (5,2)") is no longer implemented: strike the
# existing two cells and note it in a new column C cell.
$ws.Range("A7").Font.Strikethrough = $true
$ws.Range("B7").Font.Strikethrough = $true
$ws.Range("C7").Value = "not yet implemented"

# Column C needs to grow to fit the new, longer note text.
$ws.Columns("C").ColumnWidth = 17.666666666666668

# Move the cursor/selection to D5.
[void]$ws.Range("D5").Select()
